$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1 / index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 628
$ws1.Range("F3").Value = 1307
$ws1.Range("F4").Value = 1176
$ws1.Range("F5").Value = 14393
$ws1.Range("F6").Value = 16765
$ws1.Range("F7").Value = 16
$ws1.Range("F8").Value = 123
$ws1.Range("F9").Value = 30
$ws1.Range("F12").Value = 27
$ws1.Range("F18").Value = 111
$ws1.Range("F20").Value = 1278
$ws1.Range("F23").Value = 46
$ws1.Range("F24").Value = 29
$ws1.Range("F26").Value = 6845
$ws1.Range("F28").Value = 25
$ws1.Range("F29").Value = 1135
$ws1.Range("F30").Value = 15
$ws1.Range("F32").Value = 5779
$ws1.Range("F33").Value = 117
$ws1.Range("F35").Value = 196
$ws1.Range("F36").Value = 4887
$ws1.Range("F37").Value = 24

# Sheet "全部类型" (sheet4 / index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 628
$ws4.Range("F3").Value = 1307
$ws4.Range("F4").Value = 1176
$ws4.Range("F5").Value = 14393
$ws4.Range("F6").Value = 16765
$ws4.Range("F7").Value = 16
$ws4.Range("F8").Value = 123
$ws4.Range("F9").Value = 30
$ws4.Range("F12").Value = 27
$ws4.Range("F18").Value = 111
$ws4.Range("F20").Value = 1278
$ws4.Range("F24").Value = 46
$ws4.Range("F25").Value = 29
$ws4.Range("F27").Value = 6845
$ws4.Range("F29").Value = 25
$ws4.Range("F30").Value = 1135
$ws4.Range("F31").Value = 15
$ws4.Range("F35").Value = 5779
$ws4.Range("F36").Value = 117
$ws4.Range("F38").Value = 196
$ws4.Range("F39").Value = 4887
$ws4.Range("F40").Value = 24
